# Update Name of Algo
# Apply updated result values for the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.846500000000004
$ws.Range("B14").Value = 5.756399999999998
$ws.Range("B21").Value = 9.928300000000002
$ws.Range("D22").Value = -8.047200000000002
$ws.Range("B23").Value = 8.875899999999994
$ws.Range("D24").Value = -7.268400000000003
$ws.Range("B25").Value = 5.543499999999999
$ws.Range("B26").Value = 4.639300000000005
$ws.Range("D28").Value = -8.156899999999995
$ws.Range("B29").Value = 5.220500000000004
$ws.Range("D36").Value = -6.726699999999999
$ws.Range("D45").Value = -7.191599999999997
$ws.Range("D48").Value = -7.680299999999994
$ws.Range("D49").Value = -8.0114
$ws.Range("D52").Value = -7.945500000000009
$ws.Range("B53").Value = 5.935600000000001
$ws.Range("D53").Value = -8.411999999999999
$ws.Range("D54").Value = -7.939200000000008
$ws.Range("B57").Value = 4.896199999999993
$ws.Range("B59").Value = 4.8953
$ws.Range("B69").Value = 5.349399999999995
$ws.Range("D70").Value = -7.2662
$ws.Range("B79").Value = 9.611600000000005
$ws.Range("B83").Value = 5.203099999999995
$ws.Range("D86").Value = -8.210599999999998
$ws.Range("D87").Value = -8.660899999999991
$ws.Range("D89").Value = -8.316400000000002
$ws.Range("B91").Value = 5.896900000000004
$ws.Range("B93").Value = 5.597599999999997
$ws.Range("D101").Value = -8.098700000000001
$ws.Range("B103").Value = 5.929599999999996
